$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "30.212.72"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3
$ws.Range("D3").Value = "1.839.93"
$ws.Range("E3").Value = "  -1.30%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "232.86"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.02%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4675"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.75%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2719"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06277"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.73%  "

# Row 10
$ws.Range("D10").Value = "1.838.99"
$ws.Range("E10").Value = "  -1.64%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07422"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "16.13"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.08%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.939"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -2.62%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "83.61"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.73%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6192"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.90%  "

# Row 16
$ws.Range("D16").Value = "30.145.05"
$ws.Range("E16").Value = "  -0.89%  "

# Row 17
$ws.Range("E17").Value = "  +0.02%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "226.01"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007275"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.37%  "

# Row 20
$ws.Range("B20").Value = "Avalanche"
$ws.Range("C20").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.33"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.88%  "

# Row 21
$ws.Range("B21").Value = "BinanceUSD"
$ws.Range("C21").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.11%  "

# Row 22
$ws.Range("B22").Value = "Uniswap"
$ws.Range("C22").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.897"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -4.60%  "

# Row 23
$ws.Range("B23").Value = "Chainlink"
$ws.Range("C23").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.838"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.06%  "

# Row 24
$ws.Range("B24").Value = "Cosmos"
$ws.Range("C24").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.183"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.28%  "

# Row 25
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "164.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.76%  "

# Row 26
$ws.Range("B26").Value = "EthereumClassic"
$ws.Range("C26").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "17.74"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -2.23%  "

# Row 27
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.861"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.70%  "

# Row 28
$ws.Range("B28").Value = "Stellar"
$ws.Range("C28").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1032"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -2.33%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.374"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.84%  "

# Row 30
$ws.Range("B30").Value = "InternetComputer(DFINITY)"
$ws.Range("C30").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.074"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.09%  "

# Row 31
$ws.Range("B31").Value = "Filecoin"
$ws.Range("C31").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.799"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.07%  "

# Row 32
$ws.Range("B32").Value = "Hedera"
$ws.Range("C32").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.04809"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "

# Row 33
$ws.Range("B33").Value = "ARBITRUM"
$ws.Range("C33").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.140"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.52%  "

# Row 34
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7034"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.91%  "

# Row 35
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.701"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.38%  "

# Row 36
$ws.Range("B36").Value = "VeChain"
$ws.Range("C36").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.01868"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.18%  "

# Row 37
$ws.Range("B37").Value = "MXToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.649"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.38%  "

# Row 38
$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.8898"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.76%  "

# Row 39
$ws.Range("B39").Value = "Quant"
$ws.Range("C39").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "104.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.62%  "

# Row 40
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.62%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.911"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.04%  "

# Row 42
$ws.Range("B42").Value = "FraxShare"
$ws.Range("C42").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.517"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.68%  "

# Row 43
$ws.Range("B43").Value = "TheSandbox"
$ws.Range("C43").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4006"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.99%  "

# Row 44
$ws.Range("B44").Value = "Aptos"
$ws.Range("C44").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.020"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.88%  "

# Row 45
$ws.Range("B45").Value = "Algorand"
$ws.Range("C45").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1193"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.29%  "

# Row 46
$ws.Range("B46").Value = "Aave"
$ws.Range("C46").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.80"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.46%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "8.549"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.63%  "

# Row 48
$ws.Range("B48").Value = "Elrond"
$ws.Range("C48").Value = "https://coinranking.com/coin/omwkOTglq+elrond-egld"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "32.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.63%  "

# Row 49
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05513"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.27%  "

# Row 50
$ws.Range("B50").Value = "NEARProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.353"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.50%  "

# Row 51
$ws.Range("B51").Value = "Decentraland"
$ws.Range("C51").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3629"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.81%  "
